# Weekly price-sheet refresh: the oldest daily record (row 86) is replaced by
# a brand-new observation, and every following record shifts down one row to
# make room, pushing the last existing record (old row 129) into a new row 130.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D carries a date-style number format (s="2" in the original sheet).
# Remember it so the newly created row 130 can be formatted the same way.
$dateFmt = $ws.Range("D86").NumberFormat

# Shift rows 129..86 down into 130..87, working bottom-up so we never
# overwrite a source row before it has been copied.
for ($r = 130; $r -ge 87; $r--) {
    $srcRow = $r - 1
    $src = $ws.Range("A" + $srcRow + ":R" + $srcRow)
    $dst = $ws.Range("A" + $r + ":R" + $r)
    $dst.Value2 = $src.Value2
}

# The brand new row 130 has no pre-existing style for column D; give it the
# same date number format used by the rest of the column.
$ws.Range("D130").NumberFormat = $dateFmt

# Row 86 becomes the newest observation.
$ws.Range("D86").Value2 = 45119
$ws.Range("J86").Value2 = 70
$ws.Range("K86").Value2 = 7500
$ws.Range("L86").Value2 = 8000
$ws.Range("M86").Value2 = 7643
$ws.Range("P86").Value2 = 764
